# edit.ps1 -- apply the "New crime data collected" revision to the
# NYPD 71st Precinct weekly CompStat report (cs-en-us-071pct).
#
# The report header text shifts one week forward (Vol 29 No 42 -> No 43,
# week of 10/17-10/23/2022 -> 10/24-10/30/2022), and the Week-to-Date /
# 28-Day / Year-to-Date / 2-Year crime-count table (rows 14-30) is refreshed
# with the new week's figures. A handful of cells flip between a literal
# number and the sheet's "no data" placeholder text ("0" / "***.*"), which
# also carries its own right-aligned-text cell style (as opposed to the
# right-aligned-number style used for real counts/percentages) -- those are
# reproduced by copying an existing donor cell that already has the right
# combination of style + value, which carries both over in one shot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/issue number and the report week's date range ---
# These live inside multi-run shared strings, so only the affected run of
# characters is replaced to preserve the original rich-text formatting.
$ws.Range("A8").Characters(21, 2).Text = "43"
$ws.Range("C9").Characters(27, 10).Text = "10/24/2022"
$ws.Range("C9").Characters(48, 10).Text = "10/30/2022"

# --- Row 14 ---
$ws.Range("F15").Copy($ws.Range("C14"))  # -> 1
$ws.Range("G15").Copy($ws.Range("F14"))  # -> 1
$ws.Range("I14").Value = 7
$ws.Range("K14").Value = 250
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -78.125

# --- Row 15 ---
$ws.Range("M15").Value = 10
$ws.Range("N15").Value = -82.539682539682

# --- Row 16 ---
$ws.Range("C15").Copy($ws.Range("C16"))  # -> 0
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -9.090909090909
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = 41.666666666666
$ws.Range("L16").Value = 27.102803738317
$ws.Range("M16").Value = -40.088105726872
$ws.Range("N16").Value = -90.436005625879

# --- Row 17 ---
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -19.354838709677
$ws.Range("I17").Value = 288
$ws.Range("J17").Value = 253
$ws.Range("K17").Value = 13.833992094861
$ws.Range("L17").Value = 7.865168539325
$ws.Range("M17").Value = 12.062256809338
$ws.Range("N17").Value = -63.590391908976

# --- Row 18 ---
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 11
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 113
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = 48.684210526315
$ws.Range("L18").Value = -0.877192982456
$ws.Range("M18").Value = -57.196969696969
$ws.Range("N18").Value = -92.820838627700

# --- Row 19 ---
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 10.344827586206
$ws.Range("I19").Value = 289
$ws.Range("J19").Value = 222
$ws.Range("K19").Value = 30.180180180180
$ws.Range("L19").Value = 25.652173913043
$ws.Range("M19").Value = -20.165745856353
$ws.Range("N19").Value = -50.172413793103

# --- Row 20 ---
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 121
$ws.Range("J20").Value = 103
$ws.Range("K20").Value = 17.475728155339
$ws.Range("L20").Value = 30.107526881720
$ws.Range("M20").Value = -4.724409448818
$ws.Range("N20").Value = -89.357959542656

# --- Row 21 ---
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -3.703703703703
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 1.063829787234
$ws.Range("I21").Value = 965
$ws.Range("J21").Value = 770
$ws.Range("K21").Value = 25.324675324675
$ws.Range("L21").Value = 15.430622009569
$ws.Range("M21").Value = -23.046251993620
$ws.Range("N21").Value = -82.764779424897

# --- Row 22 ---
$ws.Range("M22").Value = -62.790697674418

# --- Row 24 ---
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -31.034482758620
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -36.190476190476
$ws.Range("I24").Value = 912
$ws.Range("J24").Value = 889
$ws.Range("K24").Value = 2.587176602924
$ws.Range("L24").Value = 9.747292418772
$ws.Range("M24").Value = 7.168037602820

# --- Row 25 ---
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -15.384615384615
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = -29.629629629629
$ws.Range("I25").Value = 469
$ws.Range("J25").Value = 384
$ws.Range("K25").Value = 22.135416666666
$ws.Range("L25").Value = 48.888888888888
$ws.Range("M25").Value = -33.380681818181

# --- Row 26 ---
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200

# --- Row 27 ---
$ws.Range("C15").Copy($ws.Range("C27"))  # -> 0
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 53
$ws.Range("K27").Value = -7.547169811320
$ws.Range("L27").Value = 2.083333333333

# --- Row 28 ---
$ws.Range("C28").Value = 2
$ws.Range("F15").Copy($ws.Range("D28"))  # -> 1
$ws.Range("H26").Copy($ws.Range("E28"))  # -> 100
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = -13.043478260869
$ws.Range("L28").Value = -52.380952380952
$ws.Range("M28").Value = -45.945945945945
$ws.Range("N28").Value = -77.777777777777

# --- Row 29 ---
$ws.Range("C29").Value = 2
$ws.Range("F15").Copy($ws.Range("D29"))  # -> 1
$ws.Range("H26").Copy($ws.Range("E29"))  # -> 100
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 18
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = -10
$ws.Range("L29").Value = -28
$ws.Range("M29").Value = -37.931034482758
$ws.Range("N29").Value = -79.775280898876

# --- Row 30 ---
$ws.Range("C15").Copy($ws.Range("D30"))  # -> 0
$ws.Range("E22").Copy($ws.Range("E30"))  # -> ***.*
